# Add "HTTP/2 SSL" (Mean / Std Dev) result columns (F:G) to each results sheet,
# mirroring the existing HTTP/2 columns (D:E).

$wb = $excel.ActiveWorkbook

# Per-sheet F/G values for rows 4-7 (Mean row, then Std Dev row values follow the
# same column layout as D/E): sheet name -> row -> @(F, G)
$data = @{
    "Transfer Time (s)" = @{
        4 = @(0.01700205326080322, 0.01556676824206549)
        5 = @(0.02456668138504028, 0.01307115455139281)
        6 = @(0.1406377077102661,  0.04032567721202485)
        7 = @(1.726834058761597,   0)
    }
    "Throughput (bps)" = @{
        4 = @(857495.5432740636,  377229.6488290807)
        5 = @(4819337.391183397, 1522637.096039632)
        6 = @(7967923.367428995, 1978908.963937275)
        7 = @(7469176.729335703, 0)
    }
    "Overhead Ratio" = @{
        4 = @(1.0212890625,       0)
        5 = @(1.00267578125,      0)
        6 = @(1.000751495361328,  0)
        7 = @(1.000569629669189,  0)
    }
}

foreach ($ws in $wb.Worksheets) {
    $rows = $data[$ws.Name]
    if ($rows -eq $null) { continue }

    # New column widths to match existing columns.
    $ws.Columns.Item(6).ColumnWidth = 14.17
    $ws.Columns.Item(7).ColumnWidth = 14.17

    # Row 1 header: merge F1:G1 first (like B1:C1 / D1:E1), then apply the
    # same style/value as the other headers -- merging *after* styling would
    # redraw the merged-region border and mint new style records.
    $ws.Range("F1:G1").Merge()
    $ws.Range("D1").Copy()
    $ws.Range("F1").PasteSpecial(-4122)
    $ws.Range("E1").Copy()
    $ws.Range("G1").PasteSpecial(-4122)
    $ws.Range("F1").Value = "HTTP/2 SSL"

    # Row 2 sub-header: Mean / Std Dev, same style as D2/E2.
    $ws.Range("D2").Copy()
    $ws.Range("F2").PasteSpecial(-4122)
    $ws.Range("E2").Copy()
    $ws.Range("G2").PasteSpecial(-4122)
    $ws.Range("F2").Value = "Mean"
    $ws.Range("G2").Value = "Std Dev"

    # Row 3 spacer row, same style as D3/E3 (empty).
    $ws.Range("D3").Copy()
    $ws.Range("F3").PasteSpecial(-4122)
    $ws.Range("E3").Copy()
    $ws.Range("G3").PasteSpecial(-4122)

    # Rows 4-7 data, same style as the corresponding D/E cell, with new values.
    foreach ($r in 4..7) {
        $ws.Range("D$r").Copy()
        $ws.Range("F$r").PasteSpecial(-4122)
        $ws.Range("E$r").Copy()
        $ws.Range("G$r").PasteSpecial(-4122)

        $vals = $rows[$r]
        $ws.Range("F$r").Value = $vals[0]
        $ws.Range("G$r").Value = $vals[1]
    }
}
